$d = $word.ActiveDocument

# 1) Merge the split "{{", "nombre_alumno", "}}" runs (with proofErr wrappers)
#    into a single run reading "{{nombre_alumno}}" (same rPr as before).
$d.Content.Find.Execute("{{nombre_alumno}}", $false, $false, $false, $false, $false, $true, 1, $false, "{{nombre_alumno}}", 2)

# 2) Insert the new "ASIGNATURA / trimester grades" table right after the
#    small spacer paragraph that follows the student-info table (and right
#    before the spacer paragraph that precedes the page break), matching
#    the target OOXML exactly.
$tblXml = @'
<w:tbl xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:tblPr><w:tblStyle w:val="Tablaconcuadrcula"/><w:tblpPr w:leftFromText="141" w:rightFromText="141" w:vertAnchor="text" w:horzAnchor="margin" w:tblpY="12"/><w:tblW w:w="5000" w:type="pct"/><w:tblLook w:val="04A0" w:firstRow="1" w:lastRow="0" w:firstColumn="1" w:lastColumn="0" w:noHBand="0" w:noVBand="1"/></w:tblPr><w:tblGrid><w:gridCol w:w="4719"/><w:gridCol w:w="1644"/><w:gridCol w:w="1644"/><w:gridCol w:w="1644"/><w:gridCol w:w="1895"/><w:gridCol w:w="2448"/></w:tblGrid><w:tr><w:trPr><w:trHeight w:val="737"/></w:trPr><w:tc><w:tcPr><w:tcW w:w="1728" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:t>ASIGNATURA</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="629" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:t>PRIMER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:br/><w:t>TRIMESTRE</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="629" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:t>SEGUNDO</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:br/><w:t>TRIMESTRE</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="629" w:type="pct"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:t>TERCER</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:br/><w:t>TRIMESTRE</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="1385" w:type="pct"/><w:gridSpan w:val="2"/><w:shd w:val="clear" w:color="auto" w:fill="F2F2F2" w:themeFill="background1" w:themeFillShade="F2"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/></w:rPr><w:t>NOTA TOTAL</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="5000" w:type="pct"/><w:gridSpan w:val="6"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{% for asgin in asignatura %}}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="1728" w:type="pct"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>asig.nombre_</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>asignatura}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="629" w:type="pct"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{asig.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>t1</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="629" w:type="pct"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{asig.t</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>2</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="629" w:type="pct"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{asig.t</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>3</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="469" w:type="pct"/><w:tcBorders><w:right w:val="nil"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:b/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{asig.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>nota_final</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:tc><w:tc><w:tcPr><w:tcW w:w="916" w:type="pct"/><w:tcBorders><w:top w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:left w:val="nil"/><w:bottom w:val="single" w:sz="4" w:space="0" w:color="auto"/><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:shd w:val="clear" w:color="auto" w:fill="auto"/><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:jc w:val="center"/><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>{{asig.</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>calificacion</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>}}</w:t></w:r></w:p></w:tc></w:tr><w:tr><w:tc><w:tcPr><w:tcW w:w="5000" w:type="pct"/><w:gridSpan w:val="6"/><w:tcBorders><w:right w:val="single" w:sz="4" w:space="0" w:color="auto"/></w:tcBorders><w:vAlign w:val="center"/></w:tcPr><w:p><w:pPr><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr></w:pPr><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve">{{% </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t>endfor</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Mulish" w:hAnsi="Mulish"/><w:color w:val="000000" w:themeColor="text1"/></w:rPr><w:t xml:space="preserve"> %}}</w:t></w:r></w:p></w:tc></w:tr></w:tbl>
'@

$insertionPoint = $d.Range(119, 119)
$insertionPoint.InsertXML($tblXml)
